$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 16413
$ws1.Range("G2").Value = 62
$ws1.Range("F3").Value = 354
$ws1.Range("F4").Value = 737
$ws1.Range("F6").Value = 698
$ws1.Range("F7").Value = 1739
$ws1.Range("F8").Value = 161

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 16413
$ws4.Range("G2").Value = 62
$ws4.Range("F3").Value = 354
$ws4.Range("F4").Value = 737
$ws4.Range("F8").Value = 698
$ws4.Range("F9").Value = 1739
$ws4.Range("F11").Value = 161
